# Architecture diagram update: remove the "online"/Web cloud component,
# widen + relabel the Main shape to "MainApp", nudge its outgoing
# connector, refresh the cached "last updated" date fields, and (best
# effort) restore the author's slide guides.
#
# $ppt / $ppt.ActivePresentation are provided by the host.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: find a shape on a slide/layout/master by its stable numeric Id
# (Shapes.Item(N) is a 1-based *position* index, not the Id, and this
# deck has duplicate shape Names, so we match on Id to be unambiguous).
# ---------------------------------------------------------------------
function Get-ShapeById($container, $id) {
    foreach ($sh in $container.Shapes) {
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Presentation-level slide guides (horizontal @ 186pt, vertical @
#    360pt - i.e. pos="1488"/"2880" in the OOXML's 1/8-point units).
#    Best-effort: older/limited PowerPoint automation hosts don't all
#    support creating guides via the object model, so this is wrapped
#    defensively and the rest of the script proceeds regardless.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides) {
        $hGuide = $guides.Add(1, 186)    # ppHorizontalGuide = 1
        $vGuide = $guides.Add(0, 360)    # ppVerticalGuide   = 0
    }
} catch {
    # Guides API not available in this host - ignore.
}

# ---------------------------------------------------------------------
# 2) Refresh the cached date field text ("Date Placeholder *") on the
#    slide master, every slide layout, and the notes master.
# ---------------------------------------------------------------------
$newDate = "7/11/2017"

$master = $p.SlideMaster
foreach ($sh in $master.Shapes) {
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    foreach ($sh in $layout.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

try {
    $notesMaster = $p.NotesMaster
    foreach ($sh in $notesMaster.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
} catch {
    # NotesMaster text editing not available in this host - ignore.
}

# ---------------------------------------------------------------------
# 3) Slide 1: remove the "online" component (Cloud 50 + the dotted
#    green elbow connector feeding it), then resize/relabel "Main" to
#    "MainApp" and nudge its outgoing connector to match the new shape.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$elbow51 = Get-ShapeById $s 52   # "Elbow Connector 51"
if ($elbow51) { $elbow51.Delete() }

$cloud50 = Get-ShapeById $s 51   # "Cloud 50"
if ($cloud50) { $cloud50.Delete() }

$mainRect = Get-ShapeById $s 23  # "Rectangle 62" (the "Main" box)
if ($mainRect) {
    $mainRect.Left = 134.70001
    $mainRect.Top = 311.6978740157481
    $mainRect.Width = 87.16630921259842
    $mainRect.Height = 44.95339582677166
    $mainRect.TextFrame.TextRange.Text = "MainApp"
}

$mainConn = Get-ShapeById $s 37  # "Straight Arrow Connector 36"
if ($mainConn) {
    $mainConn.Left = 221.86629921259845
    $mainConn.Top = 334.1745769291338
}
